# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn, de-de)
# describing the newly handed-off file f52b368f-d144-49b6-aed7-4c2624f1faca.md

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

$newFile   = "f52b368f-d144-49b6-aed7-4c2624f1faca"
$newMd     = "$newFile.md"
$newZhXlf  = "$newFile.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.zh-cn.xlf"
$newDeXlf  = "$newFile.bef6f3e0c91a8143aa157effb5970a9bb752cfe2.de-de.xlf"

$commitSha = "e4317ec838fcb7d85ac106cbc52263e36e686afe"
$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/$commitSha/e2e/$newMd"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c8b98565e194f28ec4d09d15d0aa71076ff781d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newZhXlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b48c7d0e90f9e5b0dc7e8bb87e3de61989584bc4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newDeXlf"

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (row 3)
# ---------------------------------------------------------------------------
$a3 = $wsOverview.Range("A3")
$a3.Value = $newMd
$wsOverview.Hyperlinks.Add($a3, $mdUrl, "", "", $newMd) | Out-Null

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$d3 = $wsOverview.Range("D3")
$d3.Value = "2016-03-23 22:38:10"
$d3.NumberFormat = $dateFormat

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn  (row 3)
# ---------------------------------------------------------------------------
$a3zh = $wsZhCn.Range("A3")
$a3zh.Value = $newMd
$wsZhCn.Hyperlinks.Add($a3zh, $mdUrl, "", "", $newMd) | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"

$d3zh = $wsZhCn.Range("D3")
$d3zh.Value = $newZhXlf
$wsZhCn.Hyperlinks.Add($d3zh, $zhXlfUrl, "", "", $newZhXlf) | Out-Null

$e3zh = $wsZhCn.Range("E3")
$e3zh.Value = "2016-03-23 22:38:05"
$e3zh.NumberFormat = $dateFormat

$h3zh = $wsZhCn.Range("H3")
$h3zh.Value = "0001-01-01 00:00:00"
$h3zh.NumberFormat = $dateFormat

$wsZhCn.Range("J3").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet 3: de-de  (row 3)
# ---------------------------------------------------------------------------
$a3de = $wsDeDe.Range("A3")
$a3de.Value = $newMd
$wsDeDe.Hyperlinks.Add($a3de, $mdUrl, "", "", $newMd) | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"

$d3de = $wsDeDe.Range("D3")
$d3de.Value = $newDeXlf
$wsDeDe.Hyperlinks.Add($d3de, $deXlfUrl, "", "", $newDeXlf) | Out-Null

$e3de = $wsDeDe.Range("E3")
$e3de.Value = "2016-03-23 22:38:10"
$e3de.NumberFormat = $dateFormat

$h3de = $wsDeDe.Range("H3")
$h3de.Value = "0001-01-01 00:00:00"
$h3de.NumberFormat = $dateFormat

$wsDeDe.Range("J3").Value = "Include"

Write-Host "Handoff report row added for $newFile"
